$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns K and L ---
$ws.Range("K1").Value = "fxppo2_accuracy_qkeras"
$ws.Range("L1").Value = "orig-fxppo2-drop_qkeras"

# Copy the header formatting (bold, centered, thin border) from the
# existing J1 header cell onto the two new header cells so they match
# the rest of the header row style.
$ws.Range("J1").Copy()
$ws.Range("K1:L1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data rows 2-21: fxppo2_accuracy_qkeras (K) / orig-fxppo2-drop_qkeras (L) ---
# Values are supplied as strings and cast to [double], since a couple of
# them are in scientific notation (e.g. "-5.551115123125783e-17") which
# this PowerShell parser cannot tokenize as a numeric literal directly.
$data = @{
    2  = @("0.494824016563147", "-0.09316770186335399")
    3  = @("0.4265010351966874", "0.1532091097308489")
    4  = @("0.3850931677018634", "0.01035196687370593")
    5  = @("0.474120082815735", "0.10351966873706")
    6  = @("0.3975155279503105", "-5.551115123125783e-17")
    7  = @("0.4244306418219462", "0.05590062111801242")
    8  = @("0.4906832298136646", "0.07246376811594207")
    9  = @("0.4120082815734989", "-5.551115123125783e-17")
    10 = @("0.5403726708074534", "0.0331262939958592")
    11 = @("0.4244306418219462", "0")
    12 = @("0.008281573498964804", "0")
    13 = @("0.4244306418219462", "-0.02277432712215316")
    14 = @("0.3975155279503105", "-0.004140786749482317")
    15 = @("0.3975155279503105", "0.1677018633540373")
    16 = @("0.5693581780538303", "0.006211180124223614")
    17 = @("0.4803312629399586", "0.09937888198757766")
    18 = @("0.008281573498964804", "0")
    19 = @("0.3975155279503105", "0.002070393374741242")
    20 = @("0.5424430641821946", "0.04140786749482406")
    21 = @("0.008281573498964804", "0")
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 11).Value = [double]$vals[0]
    $ws.Cells.Item($r, 12).Value = [double]$vals[1]
}
